# Update cryptos list (Coin / Price / Volume(1h) columns) to latest scrape values.
# Numeric-looking price strings are prefixed with a leading apostrophe so Excel
# stores them as literal text (matching the source data's inline-string format)
# instead of silently coercing them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.376.51'
$ws.Range('E2').Value = '  -2.46%  '
$ws.Range('D3').Value = '2.574.64'
$ws.Range('E3').Value = '  -3.04%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''545.81'
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('D6').Value = '''143.24'
$ws.Range('E6').Value = '  -2.16%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '''0.586'
$ws.Range('E8').Value = '  +1.91%  '
$ws.Range('D9').Value = '''6.77'
$ws.Range('E9').Value = '  +1.57%  '
$ws.Range('D10').Value = '''0.0995'
$ws.Range('E10').Value = '  -3.96%  '
$ws.Range('D11').Value = '''0.140'
$ws.Range('E11').Value = '  +3.66%  '
$ws.Range('E12').Value = '  -2.64%  '
$ws.Range('D13').Value = '3.028.44'
$ws.Range('E13').Value = '  -3.10%  '
$ws.Range('D14').Value = '58.312.85'
$ws.Range('E14').Value = '  -2.43%  '
$ws.Range('D15').Value = '''20.47'
$ws.Range('E15').Value = '  -3.76%  '
$ws.Range('D16').Value = '2.588.18'
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('D17').Value = '''0.0000130'
$ws.Range('E17').Value = '  -4.09%  '
$ws.Range('D18').Value = '''4.41'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').Value = '''332.65'
$ws.Range('E19').Value = '  -3.46%  '
$ws.Range('D20').Value = '''9.96'
$ws.Range('E20').Value = '  -4.70%  '
$ws.Range('D21').Value = '''6.04'
$ws.Range('E21').Value = '  -4.79%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '''66.43'
$ws.Range('E23').Value = '  -0.56%  '
$ws.Range('D24').Value = '''0.420'
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('D25').Value = '''0.998'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = '''0.157'
$ws.Range('E26').Value = '  -5.40%  '
$ws.Range('D27').Value = '''7.02'
$ws.Range('E27').Value = '  -4.68%  '
$ws.Range('B28').Value = 'USDe'
$ws.Range('C28').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D28').Value = '''0.999'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0731'
$ws.Range('E29').Value = '  -3.67%  '
$ws.Range('D30').Value = '''1.64'
$ws.Range('E30').Value = '  -1.03%  '
$ws.Range('D31').Value = '''154.64'
$ws.Range('E31').Value = '  +2.58%  '
$ws.Range('D32').Value = '''5.83'
$ws.Range('E32').Value = '  -0.73%  '
$ws.Range('D33').Value = '''18.76'
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('D34').Value = '''3.87'
$ws.Range('E34').Value = '  -4.29%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '''37.12'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').Value = '''0.841'
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '''1.09'
$ws.Range('E37').Value = '  -4.99%  '
$ws.Range('D38').Value = '''0.811'
$ws.Range('E38').Value = '  -4.07%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''1.41'
$ws.Range('E39').Value = '  -4.09%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '''3.55'
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = '''278.36'
$ws.Range('E41').Value = '  -4.74%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '''0.999'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '''0.590'
$ws.Range('E43').Value = '  -3.19%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').Value = '''10.63'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = '''0.0938'
$ws.Range('E45').Value = '  -1.37%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '''0.0524'
$ws.Range('E46').Value = '  -2.91%  '
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('D48').Value = '1.895.39'
$ws.Range('E48').Value = '  -4.43%  '
$ws.Range('D49').Value = '''4.37'
$ws.Range('E49').Value = '  -8.26%  '
$ws.Range('D50').Value = '''17.57'
$ws.Range('E50').Value = '  -4.84%  '
$ws.Range('D51').Value = '''111.03'
$ws.Range('E51').Value = '  +0.68%  '
